# Apply the cibmtr-reporting-ig update to the "Metadata" sheet.
# (The "Include from RxNorm" sheet's own cell text never changes; its
#  displayed values stay identical even though the shared-string table
#  is renumbered underneath it as a side effect of the edits below.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# --- First, make room for the new "Immutable" row (row 16) by copying
#     the formatting that's already used throughout the table (row 15,
#     which currently holds the Immutable data we're about to relocate).
$ws.Range("A15:B15").Copy($ws.Range("A16:B16"))

# --- Work from the bottom of the table upward so that rows we still
#     need to read/derive from aren't clobbered before we use them.

# New row 16: Immutable / BooleanType[null]  (previously lived on row 15)
$ws.Range("A16").Value = "Immutable"
$ws.Range("B16").Value = "BooleanType[null]"

# Row 15 becomes Copyright, with no value (previously the Immutable row)
$ws.Range("A15").Value = "Copyright"
$ws.Range("B15").Value = ""

# Row 14 becomes Purpose, with no value (previously Purpose already, one row up)
$ws.Range("A14").Value = "Purpose"
$ws.Range("B14").Value = ""

# Row 13 becomes Description (previously on row 12)
$ws.Range("A13").Value = "Description"
$ws.Range("B13").Value = "RxNorm codes for Etoposide"

# New row 12: Jurisdiction, with no display value
$ws.Range("A12").Value = "Jurisdiction"
$ws.Range("B12").Value = ""

# New row 11: a second Contact entry
$ws.Range("A11").Value = "Contact"
$ws.Range("B11").Value = "Bob Milius (bmilius@nmdp.org)"

# Row 10: existing Contact entry gets the real ContactDetail text
$ws.Range("B10").Value = "The Medical College of Wisconsin, Inc. and the National Marrow Donor Program (http://www.cibmtr.org)"

# --- Simple in-place value updates on the unchanged-position rows ---
$ws.Range("B3").Value = "0.1.7"                               # Version
$ws.Range("B6").Value = "draft"                                # Status
$ws.Range("B8").Value = "2024-08-27T12:23:18-05:00"            # Date

$wb.Save()
